$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("quiz")

# Update "Total" row correct/total marks values
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 70
$ws.Range("E12").Value = "70/140"
